$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -3
$ws.Range("F17").Value = 4
